$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is numeric-looking: pre-format as Text so Excel
# keeps them as literal strings (matching the source inlineStr cells) instead of
# silently converting to a Number (which would also mangle trailing zeros).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.466.49'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.677.81'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '219.08'
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").Value = '0.5334'
$ws.Range("E6").Value = '  +1.92%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.2709'
$ws.Range("E8").Value = '  +4.35%  '
$ws.Range("D9").Value = '0.06407'
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("D10").Value = '21.89'
$ws.Range("E10").Value = '  +6.28%  '
$ws.Range("D11").Value = '0.07800'
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("D12").Value = '4.517'
$ws.Range("E12").Value = '  +2.80%  '
$ws.Range("D13").Value = '1.679.73'
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '0.5593'
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("D15").Value = '0.0₅8334'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").Value = '65.74'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '26.518.32'
$ws.Range("E17").Value = '  +1.82%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '4.803'
$ws.Range("E19").Value = '  +2.72%  '
$ws.Range("D20").Value = '193.31'
$ws.Range("E20").Value = '  +2.79%  '
$ws.Range("D21").Value = '10.31'
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("D22").Value = '6.324'
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '141.13'
$ws.Range("E24").Value = '  -2.85%  '
$ws.Range("D25").Value = '0.1279'
$ws.Range("E25").Value = '  +6.23%  '
$ws.Range("D26").Value = '7.414'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '16.29'
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("D28").Value = '1.444'
$ws.Range("E28").Value = '  +4.85%  '
$ws.Range("D29").Value = '0.06297'
$ws.Range("E29").Value = '  +5.89%  '
$ws.Range("D30").Value = '1.289'
$ws.Range("E30").Value = '  +2.82%  '
$ws.Range("D31").Value = '3.608'
$ws.Range("E31").Value = '  +5.40%  '
$ws.Range("D32").Value = '3.461'
$ws.Range("E32").Value = '  +1.99%  '
$ws.Range("D33").Value = '1.700'
$ws.Range("E33").Value = '  +3.60%  '
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("D35").Value = '0.6163'
$ws.Range("E35").Value = '  +9.39%  '
$ws.Range("D37").Value = '2.782'
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01632'
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '6.142'
$ws.Range("E39").Value = '  +8.08%  '
$ws.Range("D40").Value = '1.096.12'
$ws.Range("E40").Value = '  +6.71%  '
$ws.Range("D41").Value = '0.8675'
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '100.67'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '1.823.78'
$ws.Range("E45").Value = '  +4.37%  '
$ws.Range("D46").Value = '58.57'
$ws.Range("E46").Value = '  +5.19%  '
$ws.Range("D47").Value = '8.191'
$ws.Range("E47").Value = '  +1.81%  '
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = '0.05211'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.471'
$ws.Range("E50").Value = '  +6.90%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.029'
$ws.Range("E51").Value = '  +2.29%  '
